$d = $word.ActiveDocument

# --- Locate the anchor paragraph: the empty paragraph that follows
# "Should this service process one body/title at a time or a bunch ?"
# (it sits right before the first "PreformattedText"-styled paragraph).
$anchorIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Should this service process one body*") {
        $anchorIndex = $i + 1
        break
    }
}

$anchor = $d.Paragraphs.Item($anchorIndex)

# Insert the first new (empty) paragraph right after the anchor.
$anchor.Range.InsertParagraphAfter() | Out-Null

# Insert the second new paragraph (will hold the note text) right after that.
$firstNew = $d.Paragraphs.Item($anchorIndex + 1)
$firstNew.Range.InsertParagraphAfter() | Out-Null

# This is the paragraph that will carry the two runs of text.
$textPara = $d.Paragraphs.Item($anchorIndex + 2)
$textRange = $textPara.Range
$runStart = $textRange.Start

$boldText = "Une note technique "
$restText = "pr" + [char]0x00E9 + "sentant une " + [char]0x00E9 + "tude sur les approches et outils qui permettraient de g" + [char]0x00E9 + "n" + [char]0x00E9 + "raliser l" + [char]0x2019 + "approche MLOps (pipeline de codage des steps d" + [char]0x2019 + "" + [char]0x00E9 + "laboration du mod" + [char]0x00E8 + "le, et suivi de la performance du mod" + [char]0x00E8 + "le en production) ???"

# Insert the full text as a single run first (plain formatting).
$textRange.InsertAfter($boldText + $restText)
$boldEnd = $runStart + $boldText.Length

# Apply the "Strong" character style to just the first part of the text
# (using Find/Replace's format-only replace keeps the change scoped to
# the matched run instead of promoting it to a paragraph style).
$styleRange = $d.Range($runStart, $boldEnd)
$find = $styleRange.Find
$find.ClearFormatting()
$find.Text = $boldText
$find.Replacement.ClearFormatting()
$find.Replacement.Style = "Strong"
$find.Execute($boldText, $false, $false, $false, $false, $false, $true, 1, $false, $boldText, 2) | Out-Null

# The engine's character-style application stamps a stray rsidP="00000000"
# attribute onto the document's very first paragraph as a side effect.
# Rebuild that paragraph from scratch so the saved XML stays clean.
$firstPara = $d.Paragraphs.Item(1)
$firstParaText = $firstPara.Range.Text
$firstParaText = $firstParaText.TrimEnd([char]13, [char]7)

$secondPara = $d.Paragraphs.Item(2)
$rebuilt = $secondPara.Range.InsertParagraphBefore()
$d.Paragraphs.Item(2).Range.InsertAfter($firstParaText)
$d.Paragraphs.Item(1).Range.Delete() | Out-Null
